# Apply updated "想去人数" (F column) values to the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Sheet "展览": row -> new F value
$exhibitUpdates = @{
    2  = 1912
    3  = 516
    6  = 2792
    17 = 12
    19 = 227
    21 = 16
    22 = 1
    24 = 235
    25 = 14
    26 = 74
    27 = 1780
    29 = 428
    30 = 92
    33 = 315
    34 = 459
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Sheet "全部类型": row -> new F value
$allUpdates = @{
    2  = 1912
    4  = 516
    7  = 2793
    18 = 12
    20 = 227
    22 = 16
    23 = 1
    25 = 235
    26 = 14
    27 = 74
    28 = 1780
    30 = 428
    31 = 92
    34 = 315
    35 = 459
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}

$wb.Save()
